$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "81.444.46"
$ws.Range("E2").Value = "  +4.81%  "
$ws.Range("D3").Value = "3.177.50"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "207.82"
$ws.Range("E5").Value = "  +2.49%  "
$ws.Range("D6").Value = "633.32"
$ws.Range("E6").Value = "  +0.81%  "
$ws.Range("D7").Value = "0.292"
$ws.Range("E7").Value = "  +28.31%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "0.590"
$ws.Range("E9").Value = "  +2.60%  "
$ws.Range("D10").Value = "3.174.16"
$ws.Range("E10").Value = "  +0.28%  "
$ws.Range("D11").Value = "0.589"
$ws.Range("E11").Value = "  +4.58%  "
$ws.Range("D12").Value = "0.0000261"
$ws.Range("E12").Value = "  +16.27%  "
$ws.Range("D13").Value = "0.165"
$ws.Range("E13").Value = "  +2.09%  "
$ws.Range("D14").Value = "5.35"
$ws.Range("E14").Value = "  -1.39%  "
$ws.Range("D15").Value = "3.760.35"
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("D16").Value = "31.91"
$ws.Range("E16").Value = "  +1.78%  "
$ws.Range("D17").Value = "81.297.14"
$ws.Range("E17").Value = "  +5.03%  "
$ws.Range("D18").Value = "3.178.80"
$ws.Range("E18").Value = "  +1.20%  "
$ws.Range("D19").Value = "3.23"
$ws.Range("E19").Value = "  +14.09%  "
$ws.Range("D20").Value = "14.15"
$ws.Range("E20").Value = "  -1.48%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "439.37"
$ws.Range("E21").Value = "  +2.22%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "9.18"
$ws.Range("E22").Value = "  -2.37%  "
$ws.Range("D23").Value = "5.16"
$ws.Range("E23").Value = "  +6.32%  "
$ws.Range("D24").Value = "7.13"
$ws.Range("E24").Value = "  +4.19%  "
$ws.Range("E25").Value = "  +9.49%  "
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "3.345.03"
$ws.Range("E26").Value = "  +1.13%  "
$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").Value = "11.18"
$ws.Range("E27").Value = "  +3.50%  "
$ws.Range("D28").Value = "76.88"
$ws.Range("E28").Value = "  +0.99%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("E30").Value = "  +9.81%  "
$ws.Range("D31").Value = "9.15"
$ws.Range("E31").Value = "  +3.90%  "
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("D33").Value = "562.52"
$ws.Range("E33").Value = "  +8.04%  "
$ws.Range("E34").Value = "  +2.04%  "
$ws.Range("D35").Value = "2.03"
$ws.Range("E35").Value = "  +3.34%  "
$ws.Range("D36").Value = "0.152"
$ws.Range("E36").Value = "  +12.14%  "
$ws.Range("D37").Value = "0.138"
$ws.Range("E37").Value = "  +27.34%  "
$ws.Range("D38").Value = "23.14"
$ws.Range("E38").Value = "  +3.07%  "
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("D40").Value = "0.414"
$ws.Range("E40").Value = "  +3.95%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "3.08"
$ws.Range("E41").Value = "  +21.74%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").Value = "5.94"
$ws.Range("E42").Value = "  +10.23%  "
$ws.Range("E43").Value = "  +15.99%  "
$ws.Range("D44").Value = "20.76"
$ws.Range("E44").Value = "  +3.53%  "
$ws.Range("D45").Value = "160.11"
$ws.Range("E45").Value = "  -2.17%  "
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").Value = "188.82"
$ws.Range("E47").Value = "  -3.73%  "
$ws.Range("D48").Value = "1.34"
$ws.Range("E48").Value = "  +3.90%  "
$ws.Range("D49").Value = "44.18"
$ws.Range("E49").Value = "  +3.58%  "
$ws.Range("D50").Value = "0.787"
$ws.Range("E50").Value = "  -1.28%  "
$ws.Range("D51").Value = "4.26"
$ws.Range("E51").Value = "  +4.08%  "
